{"js": "// Round the displayed CNV metric values to 3 decimal places.\n// Each old value is replaced by its rounded text (e.g. \"4.72832729461341\" -> \"4.728\").\nconst replacements = [\n  [\"4.72832729461341\", \"4.728\"],\n  [\"12.2447144210211\", \"12.245\"],\n  [\"0.0697214279801276\", \"0.070\"],\n  [\"0.0\", \"0.000\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Round the displayed CNV metric values to 3 decimal places.\n# Each old value is replaced by its rounded text (e.g. \"4.72832729461341\" -> \"4.728\").\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"4.72832729461341\"; New = \"4.728\" },\n    @{ Old = \"12.2447144210211\"; New = \"12.245\" },\n    @{ Old = \"0.0697214279801276\"; New = \"0.070\" },\n    @{ Old = \"0.0\"; New = \"0.000\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $true\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$pair.Old, $true, $true, $false, $false, $false, $true, 1, $false, [ref]$pair.New, 2) | Out-Null\n}\n"}
